# Swap the data values of rows 4 and 5 for the columns that actually
# differ between the two records (A, B, E, F, G, H, Q, R, AO).
# All other columns (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE,
# AG, AI, AT, AW, AX, AY) are identical between the two rows, so no
# swap is needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R", "AO")

foreach ($col in $columns) {
    $cellRow4 = $ws.Range("$col`4")
    $cellRow5 = $ws.Range("$col`5")

    $valueRow4 = $cellRow4.Value()
    $valueRow5 = $cellRow5.Value()

    $cellRow4.Value = $valueRow5
    $cellRow5.Value = $valueRow4
}
